$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grading points entered in column E (Points for grading) for the
# "Customer Class" and "Product Class" sections.
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 1

# Grading comment explaining the point deduction on the "toString() method"
# row of the Product Class table.
$ws.Range("F14").Value = "(-1) for not seperating two instance variables by hyphen"

# Match the formatting (wrap text + full border) already used by the
# neighbouring grading-instructions cell instead of the empty comment style.
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect where the author was last working: scrolled back to the top of
# the sheet with the new comment cell selected.
$ws.Range("F14").Select()

$wb.RecalculateFull = $true
